$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.102.16"
$ws.Range("E2").Value = "  +2.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.776.23"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "339.13"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3819"
$ws.Range("E7").Value = "  -2.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3417"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.00"
$ws.Range("E9").Value = "  -2.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.143"
$ws.Range("E10").Value = "  -4.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07368"
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.09"
$ws.Range("E12").Value = "  +5.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.003"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.380"
$ws.Range("E14").Value = "  -2.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.401"
$ws.Range("E15").Value = "  +3.55%  "
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001077"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06645"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.48"
$ws.Range("E19").Value = "  -2.70%  "
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.35"
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.387"
$ws.Range("E22").Value = "  -2.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.109.24"
$ws.Range("E23").Value = "  +2.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.05"
$ws.Range("E24").Value = "  -3.01%  "
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.466"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.69"
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.407"
$ws.Range("E28").Value = "  -4.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.25"
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.980.35"
$ws.Range("E30").Value = "  -0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "134.60"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.035"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.044"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08912"
$ws.Range("E34").Value = "  +1.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.69"
$ws.Range("E35").Value = "  -3.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02405"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6813"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06360"
$ws.Range("E38").Value = "  -2.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.287"
$ws.Range("E39").Value = "  -2.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2152"
$ws.Range("E40").Value = "  -2.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.237"
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.497"
$ws.Range("E42").Value = "  -7.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.178"
$ws.Range("E43").Value = "  -2.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.22"
$ws.Range("E44").Value = "  -2.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6255"
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.862"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.64"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.067"
$ws.Range("E49").Value = "  -3.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07502"
$ws.Range("E50").Value = "  +4.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.206"
$ws.Range("E51").Value = "  +3.98%  "
